$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that keeps its original (untouched) style, used to restore
# formatting on cells where we must force a text/string value that looks
# like a number (Excel would otherwise auto-convert it to a numeric cell).
$fmtSrc = $ws.Cells.Item(1, 1)

function Set-TextValue($cell, $text) {
    # Determine whether Excel would interpret this text as a number.
    $isNumeric = $text -match '^-?[0-9]+(\.[0-9]+)?$'

    if ($isNumeric) {
        # Force text interpretation with a leading apostrophe, then restore
        # the original cell formatting (the apostrophe trick changes style).
        $cell.Value = "'" + $text
        $fmtSrc.Copy()
        $cell.PasteSpecial(-4122)   # xlPasteFormats
    } else {
        $cell.Value = $text
    }
}

# The edit rotates columns B, C, D for every row (1-9):
#   new B = old C
#   new C = old D
#   new D = old B
# For data rows (2-9) the value that lands in column C is numeric (a count),
# so it is written as a real number instead of text.
for ($r = 1; $r -le 9; $r++) {
    $oldB = $ws.Cells.Item($r, 2).Value2
    $oldC = $ws.Cells.Item($r, 3).Value2
    $oldD = $ws.Cells.Item($r, 4).Value2

    Set-TextValue $ws.Cells.Item($r, 2) $oldC

    if ($r -eq 1) {
        Set-TextValue $ws.Cells.Item($r, 3) $oldD
    } else {
        $ws.Cells.Item($r, 3).Value = $oldD
    }

    Set-TextValue $ws.Cells.Item($r, 4) $oldB
}

# Row 8's word list (now in column D) also had its internal comma-separated
# order edited as part of this change.
Set-TextValue $ws.Cells.Item(8, 4) "phòng chống, Tư vấn, Phòng chống, du lịch, xây lắp, vận chuyển, Vận tải, dự phòng, vận tải"
